$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - standalone (non-shared) formulas
$ws.Range("E3").Formula = "=C3*D3"
$ws.Range("G3").Value = 1.0900000000000001
$ws.Range("H3").Formula = "=C3*G3"

# Row 4 venue price overrides
$ws.Range("G4").Value = 1.091
$ws.Range("G5").Value = 1.0920000000000001
$ws.Range("G6").Value = 1.093
$ws.Range("G7").Value = 1.0940000000000001

# Rows 4:7 - shared formulas (assigning one formula to a multi-cell range
# produces Excel's shared-formula <f t="shared"> encoding)
$ws.Range("E4:E7").Formula = "=C4*D4"
$ws.Range("H4:H7").Formula = "=C4*G4"

# Row 9 - sums
$ws.Range("C9").Formula = "=SUM(C3:C8)"
$ws.Range("E9").Formula = "=SUM(E3:E8)"
$ws.Range("H9").Formula = "=SUM(H3:H8)"

# Row 11 - VWAP
$ws.Range("D11").Formula = "=E9/C9"
$ws.Range("G11").Formula = "=H9/C9"

$ws.Range("G12").Select()
